# Weekly update: insert the newest price record (week of 2023-12-20) at the
# top of the data (row 3, right after the header + the oldest/benchmark row),
# pushing the existing rows 3-9 down to rows 4-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 3; Excel shifts rows 3:9 down to 4:10 and
# carries the row-above formatting (keeps the date style on column D).
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with this week's data.
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = 45280
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100101
$ws.Range("H3").Value = "Berries"
$ws.Range("I3").Value = 100101004
$ws.Range("J3").Value = "Frambuesa"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 8000
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 8000
$ws.Range("Q3").Value = "$/bandeja 2 kilos"
$ws.Range("R3").Value = "Región de Ñuble"
$ws.Range("S3").Value = 4000
$ws.Range("T3").Value = 2
